# Updated cryptos list values (Price / Volume(1h)) per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.100.92"
$ws.Range("E2").Value = "'  -3.11%  "
$ws.Range("D3").Value = "'1.870.33"
$ws.Range("E3").Value = "'  -2.10%  "
$ws.Range("E4").Value = "'  +0.38%  "
$ws.Range("E5").Value = "'  -1.92%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.23%  "
$ws.Range("D7").Value = "'0.5055"
$ws.Range("E7").Value = "'  +1.18%  "
$ws.Range("D8").Value = "'0.3753"
$ws.Range("E8").Value = "'  -1.70%  "
$ws.Range("D9").Value = "'0.07151"
$ws.Range("E9").Value = "'  -2.09%  "
$ws.Range("D10").Value = "'0.8882"
$ws.Range("E10").Value = "'  -2.76%  "
$ws.Range("E11").Value = "'  -2.81%  "
$ws.Range("D12").Value = "'0.07558"
$ws.Range("E12").Value = "'  -1.65%  "
$ws.Range("D13").Value = "'1.860.36"
$ws.Range("E13").Value = "'  -1.44%  "
$ws.Range("D14").Value = "'5.323"
$ws.Range("E14").Value = "'  -3.48%  "
$ws.Range("D15").Value = "'89.30"
$ws.Range("E15").Value = "'  -3.70%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "'  +0.41%  "
$ws.Range("D17").Value = "'0.000008479"
$ws.Range("E17").Value = "'  -3.05%  "
$ws.Range("D18").Value = "'14.12"
$ws.Range("E18").Value = "'  -3.74%  "
$ws.Range("E19").Value = "'  +0.16%  "
$ws.Range("D20").Value = "'27.142.56"
$ws.Range("E20").Value = "'  -3.01%  "
$ws.Range("D21").Value = "'5.077"
$ws.Range("E21").Value = "'  -2.06%  "
$ws.Range("D22").Value = "'2.102.09"
$ws.Range("E22").Value = "'  -1.07%  "
$ws.Range("D23").Value = "'10.56"
$ws.Range("E23").Value = "'  -2.65%  "
$ws.Range("D24").Value = "'6.484"
$ws.Range("E24").Value = "'  -1.73%  "
$ws.Range("D25").Value = "'150.90"
$ws.Range("D26").Value = "'1.841"
$ws.Range("E26").Value = "'  -0.46%  "
$ws.Range("E27").Value = "'  -2.28%  "
$ws.Range("D28").Value = "'2.099"
$ws.Range("E28").Value = "'  -5.58%  "
$ws.Range("D29").Value = "'112.69"
$ws.Range("E29").Value = "'  -2.29%  "
$ws.Range("D30").Value = "'4.757"
$ws.Range("E30").Value = "'  -3.12%  "
$ws.Range("E31").Value = "'  -3.64%  "
$ws.Range("D32").Value = "'0.09023"
$ws.Range("E32").Value = "'  -0.02%  "
$ws.Range("D33").Value = "'0.05130"
$ws.Range("E33").Value = "'  -2.95%  "
$ws.Range("D34").Value = "'3.098"
$ws.Range("E34").Value = "'  -3.39%  "
$ws.Range("D35").Value = "'0.7419"
$ws.Range("E35").Value = "'  -3.99%  "
$ws.Range("D36").Value = "'1.159"
$ws.Range("E36").Value = "'  -6.12%  "
$ws.Range("D37").Value = "'0.02036"
$ws.Range("E37").Value = "'  -2.66%  "
$ws.Range("D38").Value = "'2.526"
$ws.Range("E38").Value = "'  -1.66%  "
$ws.Range("E39").Value = "'  -0.77%  "
$ws.Range("D40").Value = "'1.078"
$ws.Range("E40").Value = "'  -1.52%  "
$ws.Range("D41").Value = "'0.5366"
$ws.Range("E41").Value = "'  -3.70%  "
$ws.Range("D42").Value = "'6.589"
$ws.Range("E42").Value = "'  -4.36%  "
$ws.Range("D43").Value = "'115.55"
$ws.Range("E43").Value = "'  +2.41%  "
$ws.Range("D44").Value = "'8.420"
$ws.Range("E44").Value = "'  -1.16%  "
$ws.Range("D45").Value = "'0.1472"
$ws.Range("E45").Value = "'  -3.32%  "
$ws.Range("D46").Value = "'0.4640"
$ws.Range("E46").Value = "'  -4.20%  "
$ws.Range("D48").Value = "'9.986"
$ws.Range("E48").Value = "'  -6.05%  "
$ws.Range("D49").Value = "'1.566"
$ws.Range("E49").Value = "'  -4.42%  "
$ws.Range("D50").Value = "'64.60"
$ws.Range("E50").Value = "'  -4.33%  "
$ws.Range("D51").Value = "'36.54"
$ws.Range("E51").Value = "'  -1.85%  "
